$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to be treated as text so numeric-looking strings
# like "43.971.66" are not coerced into floating point numbers.
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '43.971.66'
$ws.Range("E2").Value = '  +1.91%  '
$ws.Range("D3").Value = '2.263.48'
$ws.Range("E3").Value = '  +1.33%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = '318.56'
$ws.Range("E5").Value = '  -0.72%  '
$ws.Range("D6").Value = '101.65'
$ws.Range("E6").Value = '  +1.27%  '
$ws.Range("D7").Value = '0.576'
$ws.Range("E7").Value = '  -1.17%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").Value = '0.556'
$ws.Range("E9").Value = '  -1.41%  '
$ws.Range("E10").Value = '  +0.28%  '
$ws.Range("E11").Value = '  +0.14%  '
$ws.Range("D12").Value = '7.65'
$ws.Range("E12").Value = '  -0.64%  '
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("D14").Value = '2.608.81'
$ws.Range("E14").Value = '  +1.38%  '
$ws.Range("D15").Value = '0.864'
$ws.Range("E15").Value = '  -0.56%  '
$ws.Range("D16").Value = '14.48'
$ws.Range("E16").Value = '  +0.08%  '
$ws.Range("D17").Value = '2.269.18'
$ws.Range("E17").Value = '  +1.57%  '
$ws.Range("D18").Value = '43.886.14'
$ws.Range("E18").Value = '  +1.88%  '
$ws.Range("D19").Value = '13.38'
$ws.Range("E19").Value = '  -9.35%  '
$ws.Range("E20").Value = '  +1.59%  '
$ws.Range("D21").Value = '6.57'
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").Value = '65.70'
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  -1.63%  '
$ws.Range("D24").Value = '235.48'
$ws.Range("E24").Value = '  -1.33%  '
$ws.Range("E25").Value = '  -2.51%  '
$ws.Range("E26").Value = '  +0.19%  '
$ws.Range("E27").Value = '  +0.95%  '
$ws.Range("E28").Value = '  -3.66%  '
$ws.Range("D29").Value = '37.22'
$ws.Range("E29").Value = '  +2.82%  '
$ws.Range("E30").Value = '  -2.01%  '
$ws.Range("D31").Value = '20.24'
$ws.Range("E31").Value = '  -0.90%  '
$ws.Range("D32").Value = '158.37'
$ws.Range("E32").Value = '  +2.65%  '
$ws.Range("E33").Value = '  -2.70%  '
$ws.Range("E34").Value = '  +0.75%  '
$ws.Range("D35").Value = '0.115'
$ws.Range("E35").Value = '  +10.62%  '
$ws.Range("E36").Value = '  +0.70%  '
$ws.Range("D37").Value = '3.08'
$ws.Range("E37").Value = '  -4.03%  '
$ws.Range("E38").Value = '  -2.47%  '
$ws.Range("D39").Value = '16.20'
$ws.Range("E39").Value = '  +17.95%  '
$ws.Range("D40").Value = '3.72'
$ws.Range("E40").Value = '  +1.50%  '
$ws.Range("E41").Value = '  -5.49%  '
$ws.Range("E42").Value = '  -3.04%  '
$ws.Range("E43").Value = '  +0.17%  '
$ws.Range("D44").Value = '1.800.65'
$ws.Range("E44").Value = '  +3.69%  '
$ws.Range("D45").Value = '76.32'
$ws.Range("E45").Value = '  +0.88%  '
$ws.Range("E46").Value = '  -3.41%  '
$ws.Range("D47").Value = '82.53'
$ws.Range("E47").Value = '  -4.26%  '
$ws.Range("E48").Value = '  -2.03%  '
$ws.Range("D49").Value = '105.07'
$ws.Range("E49").Value = '  +1.66%  '
$ws.Range("D50").Value = '58.47'
$ws.Range("E50").Value = '  -0.84%  '
$ws.Range("E51").Value = '  +4.53%  '

# Restore original (unstyled) formatting on column D now that the
# text values are locked in, so no stray style index is introduced.
$priceRange.ClearFormats()
